# =================================================================
# Commit: [ADDITIONAL SCRAPING] added code to scrape more data about
# a player's batting performance in a match, also updated the excel sheets
# =================================================================
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Insert the new "Player Info" sheet. Worksheets.Add() with no
#    args inserts immediately before the active sheet -- since
#    "ODI Batting" is active/first, this lands Player Info at the
#    very front, matching the target tab order.
# -----------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")

# -----------------------------------------------------------------
# 2) Append the new "ODI Batting Extra" sheet after "ODI Bowling".
# -----------------------------------------------------------------
$battingExtra = $wb.Worksheets.Add($null, $bowlingWs)
$battingExtra.Name = "ODI Batting Extra"

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
    $range.Borders.LineStyle = 1
}

# -----------------------------------------------------------------
# Player Info: headers + single player row
# -----------------------------------------------------------------
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
Set-HeaderStyle($playerInfo.Range("A1:D1"))

$playerInfo.Range("A2").Value = "'3707"
$playerInfo.Range("B2").Value = "Alasdair C Evans"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# -----------------------------------------------------------------
# ODI Batting: rename D1 header, rewrite D column match-card links
# as bare match codes (text), and clear the previously-blank
# INNING_NUMBER (B) placeholder cells for matches with no value.
# -----------------------------------------------------------------
$battingWs.Range("D1").Value = "MATCH_CODE"

$battingWs.Range("D2").Value = "'2979"
$battingWs.Range("D3").Value = "'2980"
$battingWs.Range("D4").Value = "'3428"
$battingWs.Range("D5").Value = "'3637"
$battingWs.Range("D6").Value = "'3675"
$battingWs.Range("D7").Value = "'3676"
$battingWs.Range("D8").Value = "'3725"
$battingWs.Range("D9").Value = "'3733"
$battingWs.Range("D10").Value = "'3761"
$battingWs.Range("D11").Value = "'3764"
$battingWs.Range("D12").Value = "'3774"
$battingWs.Range("D13").Value = "'3782"
$battingWs.Range("D14").Value = "'3880"
$battingWs.Range("D15").Value = "'3912"
$battingWs.Range("D16").Value = "'3913"
$battingWs.Range("D17").Value = "'3919"
$battingWs.Range("D18").Value = "'3920"
$battingWs.Range("D19").Value = "'3933"
$battingWs.Range("D20").Value = "'3934"
$battingWs.Range("D21").Value = "'4048"
$battingWs.Range("D22").Value = "'4049"
$battingWs.Range("D23").Value = "'4077"
$battingWs.Range("D24").Value = "'4078"
$battingWs.Range("D25").Value = "'4090"
$battingWs.Range("D26").Value = "'4091"
$battingWs.Range("D27").Value = "'4111"
$battingWs.Range("D28").Value = "'4113"
$battingWs.Range("D29").Value = "'4118"
$battingWs.Range("D30").Value = "'4142"
$battingWs.Range("D31").Value = "'4161"
$battingWs.Range("D32").Value = "'4165"
$battingWs.Range("D33").Value = "'4290"
$battingWs.Range("D34").Value = "'4302"
$battingWs.Range("D35").Value = "'4363"
$battingWs.Range("D36").Value = "'4381"
$battingWs.Range("D37").Value = "'4384"
$battingWs.Range("D38").Value = "'4461"
$battingWs.Range("D39").Value = "'4462"
$battingWs.Range("D40").Value = "'4510"
$battingWs.Range("D41").Value = "'4513"
$battingWs.Range("D42").Value = "'4515"

$emptyInningRows = 2,3,4,7,9,15,17,18,19,20,21,25,26,30,31,32,33,37,38,39,40,41,42
foreach ($r in $emptyInningRows) {
    $battingWs.Cells.Item($r, 2).ClearContents()
}

# -----------------------------------------------------------------
# ODI Bowling: rename B1 header, rewrite B column match-card links
# as bare match codes (text).
# -----------------------------------------------------------------
$bowlingWs.Range("B1").Value = "MATCH_CODE"

$bowlingWs.Range("B2").Value = "'2979"
$bowlingWs.Range("B3").Value = "'2980"
$bowlingWs.Range("B4").Value = "'3428"
$bowlingWs.Range("B5").Value = "'3637"
$bowlingWs.Range("B6").Value = "'3675"
$bowlingWs.Range("B7").Value = "'3676"
$bowlingWs.Range("B8").Value = "'3725"
$bowlingWs.Range("B9").Value = "'3761"
$bowlingWs.Range("B10").Value = "'3764"
$bowlingWs.Range("B11").Value = "'3774"
$bowlingWs.Range("B12").Value = "'3782"
$bowlingWs.Range("B13").Value = "'3880"
$bowlingWs.Range("B14").Value = "'3912"
$bowlingWs.Range("B15").Value = "'3913"
$bowlingWs.Range("B16").Value = "'3919"
$bowlingWs.Range("B17").Value = "'3920"
$bowlingWs.Range("B18").Value = "'3933"
$bowlingWs.Range("B19").Value = "'3934"
$bowlingWs.Range("B20").Value = "'4048"
$bowlingWs.Range("B21").Value = "'4049"
$bowlingWs.Range("B22").Value = "'4077"
$bowlingWs.Range("B23").Value = "'4078"
$bowlingWs.Range("B24").Value = "'4090"
$bowlingWs.Range("B25").Value = "'4091"
$bowlingWs.Range("B26").Value = "'4111"
$bowlingWs.Range("B27").Value = "'4113"
$bowlingWs.Range("B28").Value = "'4118"
$bowlingWs.Range("B29").Value = "'4142"
$bowlingWs.Range("B30").Value = "'4161"
$bowlingWs.Range("B31").Value = "'4165"
$bowlingWs.Range("B32").Value = "'4290"
$bowlingWs.Range("B33").Value = "'4302"
$bowlingWs.Range("B34").Value = "'4363"
$bowlingWs.Range("B35").Value = "'4381"
$bowlingWs.Range("B36").Value = "'4384"
$bowlingWs.Range("B37").Value = "'4461"
$bowlingWs.Range("B38").Value = "'4462"
$bowlingWs.Range("B39").Value = "'4510"
$bowlingWs.Range("B40").Value = "'4513"

# -----------------------------------------------------------------
# ODI Batting Extra: headers + 20 data rows
# -----------------------------------------------------------------
$battingExtra.Range("A1").Value = "MATCH_CODE"
$battingExtra.Range("B1").Value = "BATTING_POSITION"
$battingExtra.Range("C1").Value = "NUM_4"
$battingExtra.Range("D1").Value = "NUM_6"
$battingExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$battingExtra.Range("F1").Value = "MAN_OF_MATCH"
Set-HeaderStyle($battingExtra.Range("A1:F1"))

$battingExtra.Range("A2").Value = "'4077"
$battingExtra.Range("B2").Value = 10
$battingExtra.Range("C2").Value = "'0"
$battingExtra.Range("D2").Value = "'0"
$battingExtra.Range("F2").Value = "NO"

$battingExtra.Range("A3").Value = "'4078"
$battingExtra.Range("B3").Value = 11
$battingExtra.Range("C3").Value = "'0"
$battingExtra.Range("D3").Value = "'0"
$battingExtra.Range("F3").Value = "NO"

$battingExtra.Range("A4").Value = "'4090"
$battingExtra.Range("F4").Value = "NO"

$battingExtra.Range("A5").Value = "'4091"
$battingExtra.Range("F5").Value = "NO"

$battingExtra.Range("A6").Value = "'4111"
$battingExtra.Range("F6").Value = "NO"

$battingExtra.Range("A7").Value = "'4113"
$battingExtra.Range("B7").Value = 10
$battingExtra.Range("C7").Value = "'0"
$battingExtra.Range("D7").Value = "'1"
$battingExtra.Range("E7").Value = "4.23%"
$battingExtra.Range("F7").Value = "NO"

$battingExtra.Range("A8").Value = "'4118"
$battingExtra.Range("F8").Value = "NO"

$battingExtra.Range("A9").Value = "'4142"
$battingExtra.Range("B9").Value = 10
$battingExtra.Range("F9").Value = "NO"

$battingExtra.Range("A10").Value = "'4161"
$battingExtra.Range("B10").Value = 10
$battingExtra.Range("F10").Value = "NO"

$battingExtra.Range("A11").Value = "'4165"
$battingExtra.Range("B11").Value = 10
$battingExtra.Range("F11").Value = "NO"

$battingExtra.Range("A12").Value = "'4290"
$battingExtra.Range("B12").Value = 10
$battingExtra.Range("F12").Value = "NO"

$battingExtra.Range("A13").Value = "'4302"
$battingExtra.Range("F13").Value = "NO"

$battingExtra.Range("A14").Value = "'4363"
$battingExtra.Range("B14").Value = 10
$battingExtra.Range("C14").Value = "'0"
$battingExtra.Range("D14").Value = "'0"
$battingExtra.Range("E14").Value = "2.38%"
$battingExtra.Range("F14").Value = "NO"

$battingExtra.Range("A15").Value = "'4381"
$battingExtra.Range("F15").Value = "NO"

$battingExtra.Range("A16").Value = "'4384"
$battingExtra.Range("B16").Value = 10
$battingExtra.Range("F16").Value = "NO"

$battingExtra.Range("A17").Value = "'4461"
$battingExtra.Range("F17").Value = "NO"

$battingExtra.Range("A18").Value = "'4462"
$battingExtra.Range("F18").Value = "NO"

$battingExtra.Range("A19").Value = "'4510"
$battingExtra.Range("F19").Value = "NO"

$battingExtra.Range("A20").Value = "'4513"
$battingExtra.Range("F20").Value = "NO"

$battingExtra.Range("A21").Value = "'4515"
$battingExtra.Range("F21").Value = "NO"

# -----------------------------------------------------------------
# Restore the active tab to the first sheet (Player Info), matching
# the original workbook's activeTab="0" view state.
# -----------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
